$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Itemized billing / quote status: change C8 from OPEN to WONTFIX
$ws.Range("C8").Value = "WONTFIX"

# Update the saved cursor/selection position on the sheet
$ws.Activate()
$ws.Range("A17").Select()
